{"js": "// Update the date line.\nconst dateResults = context.document.body.search(\"2025-02-17 Monday\", { matchCase: true, matchWholeWord: false });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2025-02-18 Tuesday\", Word.InsertLocation.replace);\n}\n\n// Update every arithmetic expression cell in the practice table, in place,\n// preserving each cell's existing run formatting.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = [\n  [\"20+79=\", \"31+20=\", \"60+4=\", \"93-13=\", \"72-17=\"],\n  [\"24-23=\", \"59-15=\", \"90-35=\", \"25-8=\", \"96-4=\"],\n  [\"42+54=\", \"41-27=\", \"90+8=\", \"70-23=\", \"22+4=\"],\n  [\"8+66=\", \"72-23=\", \"73-27=\", \"91+1=\", \"93-50=\"],\n  [\"1+69=\", \"48+7=\", \"7+50=\", \"51+14=\", \"46+10=\"],\n  [\"18-2=\", \"97-81=\", \"49-6=\", \"24+1=\", \"76+18=\"],\n  [\"55-1=\", \"28+63=\", \"24+15=\", \"76+1=\", \"10+20=\"],\n  [\"74-48=\", \"24+7=\", \"18+70=\", \"29+35=\", \"12+70=\"],\n  [\"46+8=\", \"68-36=\", \"21-9=\", \"51-49=\", \"29-5=\"],\n  [\"41-18=\", \"67+5=\", \"65-21=\", \"59+12=\", \"91-53=\"],\n  [\"63-56=\", \"66-3=\", \"71+16=\", \"57-31=\", \"16-4=\"],\n  [\"38+59=\", \"10+40=\", \"98-36=\", \"59+29=\", \"53+8=\"],\n  [\"52-5=\", \"6+59=\", \"5-3=\", \"20+17=\", \"99-71=\"],\n  [\"63+36=\", \"74-49=\", \"98-2=\", \"52+32=\", \"36+54=\"],\n  [\"46-29=\", \"3+96=\", \"52-1=\", \"46+5=\", \"66-57=\"],\n  [\"35-13=\", \"78-22=\", \"63-20=\", \"99-65=\", \"38-20=\"],\n  [\"55+22=\", \"98-21=\", \"52-42=\", \"34-23=\", \"10+16=\"],\n  [\"56-6=\", \"34-1=\", \"48-36=\", \"40+57=\", \"33-30=\"],\n  [\"77+17=\", \"43-39=\", \"53-9=\", \"90-19=\", \"48+15=\"],\n  [\"37+42=\", \"44+40=\", \"73-31=\", \"62-19=\", \"81-7=\"]\n];\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date line.\n$d.Content.Find.Execute(\"2025-02-17 Monday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2025-02-18 Tuesday\", 2)\n\n# Update every arithmetic expression cell in the practice table, in place,\n# preserving each cell's existing run formatting.\n$newValues = @(\n    @(\"20+79=\", \"31+20=\", \"60+4=\", \"93-13=\", \"72-17=\"),\n    @(\"24-23=\", \"59-15=\", \"90-35=\", \"25-8=\", \"96-4=\"),\n    @(\"42+54=\", \"41-27=\", \"90+8=\", \"70-23=\", \"22+4=\"),\n    @(\"8+66=\", \"72-23=\", \"73-27=\", \"91+1=\", \"93-50=\"),\n    @(\"1+69=\", \"48+7=\", \"7+50=\", \"51+14=\", \"46+10=\"),\n    @(\"18-2=\", \"97-81=\", \"49-6=\", \"24+1=\", \"76+18=\"),\n    @(\"55-1=\", \"28+63=\", \"24+15=\", \"76+1=\", \"10+20=\"),\n    @(\"74-48=\", \"24+7=\", \"18+70=\", \"29+35=\", \"12+70=\"),\n    @(\"46+8=\", \"68-36=\", \"21-9=\", \"51-49=\", \"29-5=\"),\n    @(\"41-18=\", \"67+5=\", \"65-21=\", \"59+12=\", \"91-53=\"),\n    @(\"63-56=\", \"66-3=\", \"71+16=\", \"57-31=\", \"16-4=\"),\n    @(\"38+59=\", \"10+40=\", \"98-36=\", \"59+29=\", \"53+8=\"),\n    @(\"52-5=\", \"6+59=\", \"5-3=\", \"20+17=\", \"99-71=\"),\n    @(\"63+36=\", \"74-49=\", \"98-2=\", \"52+32=\", \"36+54=\"),\n    @(\"46-29=\", \"3+96=\", \"52-1=\", \"46+5=\", \"66-57=\"),\n    @(\"35-13=\", \"78-22=\", \"63-20=\", \"99-65=\", \"38-20=\"),\n    @(\"55+22=\", \"98-21=\", \"52-42=\", \"34-23=\", \"10+16=\"),\n    @(\"56-6=\", \"34-1=\", \"48-36=\", \"40+57=\", \"33-30=\"),\n    @(\"77+17=\", \"43-39=\", \"53-9=\", \"90-19=\", \"48+15=\"),\n    @(\"37+42=\", \"44+40=\", \"73-31=\", \"62-19=\", \"81-7=\")\n)\n\n$table = $d.Tables.Item(1)\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $row = $newValues[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $table.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
